$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-369) holds a "Förändrad" date serial that was bumped
# from 45177 (2023-09-08) to 45178 (2023-09-09) for every data row.
$ws.Range("C2:C369").Value = 45178
